$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows with refreshed market data ---
$ws.Range("D2").Value = "40.435.84"
$ws.Range("E2").Value = "  -2.77%  "
$ws.Range("D3").Value = "2.364.69"
$ws.Range("E3").Value = "  -4.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.41"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.98"
$ws.Range("E6").Value = "  -6.75%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.527"
$ws.Range("E7").Value = "  -4.36%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  -3.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0837"
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.39"
$ws.Range("E11").Value = "  -7.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.109"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "2.730.73"
$ws.Range("E13").Value = "  -4.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.52"
$ws.Range("E14").Value = "  -5.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.03"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("D16").Value = "2.337.23"
$ws.Range("E16").Value = "  -4.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.755"
$ws.Range("E17").Value = "  -4.68%  "
$ws.Range("D18").Value = "40.405.86"
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("D19").Value = "0.0₃0907"
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.11"
$ws.Range("E20").Value = "  -5.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.33"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.69"
$ws.Range("E22").Value = "  -4.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.13"
$ws.Range("E24").Value = "  -5.90%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -8.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.70"
$ws.Range("E27").Value = "  -4.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.26"
$ws.Range("E29").Value = "  -4.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.01"
$ws.Range("E30").Value = "  -6.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "154.50"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.19"
$ws.Range("E33").Value = "  -4.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0725"
$ws.Range("E34").Value = "  -4.72%  "
$ws.Range("E35").Value = "  -5.57%  "
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.05"
$ws.Range("E37").Value = "  -6.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.69"
$ws.Range("E40").Value = "  -8.34%  "
$ws.Range("E41").Value = "  -4.35%  "
$ws.Range("E42").Value = "  -5.93%  "
$ws.Range("D43").Value = "1.959.02"
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0268"
$ws.Range("E44").Value = "  -5.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.72"
$ws.Range("E45").Value = "  -5.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.37"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").Value = "  -8.84%  "
$ws.Range("D48").Value = "2.599.57"
$ws.Range("E48").Value = "  -3.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.69"
$ws.Range("E49").Value = "  -4.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.92"
$ws.Range("E50").Value = "  -5.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.94"
$ws.Range("E51").Value = "  -4.37%  "

# --- Rows 38/39 swap order: Kaspa now ranks above LidoDAOToken ---
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.100"
$ws.Range("E38").Value = "  -3.80%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.76"
$ws.Range("E39").Value = "  -4.72%  "
